$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "La mujer zorro y el doctor Shimamura"
$ws.Range("B2").Value = "Christine Wunnicke"
$ws.Range("C2").Value = "Impedimenta"
